$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates (row -> new text)
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "202"
    5  = "0.00002"
    6  = "0.00012"
    7  = "0.00005"
    9  = "0.00012"
    10 = "0.00012"
    11 = "0.00012"
    12 = "0.00671"
}

foreach ($row in $updates.Keys) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Text = $updates[$row]
}

# Rows that previously held multiple tab-separated values collapse into a
# single value (matching what rows 1-3 originally contained).
$t.Cell(44, 1).Range.Text = "100"
$t.Cell(45, 1).Range.Text = "0.01"
$t.Cell(46, 1).Range.Text = "386"
